# Daily attendance processing - 2025-12-24 10:33:37
# Applies the attendance-recording update for General Surgery B1 group sessions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Class Statistics summary block (K4:L10)
#    Recorded/Missing session counts + coverage/average-attendance %
# -----------------------------------------------------------------
$ws.Range("L6").Value = 165     # Recorded Sessions: 159 -> 165
$ws.Range("L7").Value = 3       # Missing Sessions:     9 -> 3

$ws.Range("L9").Value = "'51.9%"    # Coverage %: 50.0% -> 51.9%
$ws.Range("L10").Value = "'75.1%"   # Average Attendance %: 75.2% -> 75.1%

# restore the untouched "clean" text style (no quote-prefix) by copying
# the format from the neighboring label cells, which keeps the same
# cellXfs index as before the edit
$ws.Range("K9:K10").Copy()
$ws.Range("L9:L10").PasteSpecial(-4122)

# -----------------------------------------------------------------
# 2) "Recorded By" column (G) — System was re-ordered ahead of the
#    grader's e-mail address for every already-recorded session row.
# -----------------------------------------------------------------
$recordedByRows = @(8,9,10,34,35,36,60,61,62,86,87,88,112,113,114,138,139,140,164,167,170,191,194,197,218,221,224,245,248,251,272,275,278,299,302,305)
foreach ($r in $recordedByRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# -----------------------------------------------------------------
# 3) Group Statistics block (rows 21-26: B1D1, B1D2, B1E1, B1E2, B1F1, B1F2)
#    Recorded/Missing counts + coverage/avg-attendance % shift as the
#    corresponding 24/12/2025 B1 sessions moved from Not Recorded -> Recorded
# -----------------------------------------------------------------
$ws.Range("O21").Value = 14
$ws.Range("P21").Value = 0
$ws.Range("R21").Value = "'51.9%"
$ws.Range("S21").Value = "'77.6%"

$ws.Range("O22").Value = 14
$ws.Range("P22").Value = 0
$ws.Range("R22").Value = "'51.9%"
$ws.Range("S22").Value = "'76.0%"

$ws.Range("O23").Value = 14
$ws.Range("P23").Value = 0
$ws.Range("R23").Value = "'51.9%"
$ws.Range("S23").Value = "'80.6%"

$ws.Range("O24").Value = 13
$ws.Range("P24").Value = 1
$ws.Range("R24").Value = "'48.1%"
$ws.Range("S24").Value = "'70.9%"

$ws.Range("O25").Value = 14
$ws.Range("P25").Value = 0
$ws.Range("R25").Value = "'51.9%"
$ws.Range("S25").Value = "'69.2%"

$ws.Range("O26").Value = 14
$ws.Range("P26").Value = 0
$ws.Range("R26").Value = "'51.9%"
$ws.Range("S26").Value = "'62.8%"

# restore the pristine (non quote-prefixed) percentage-text style by
# copying formatting from the same-row "Year" label cell (column K),
# which already carries the correct style index
$ws.Range("K21:K26").Copy()
$ws.Range("R21:R26").PasteSpecial(-4122)
$ws.Range("S21:S26").PasteSpecial(-4122)

# -----------------------------------------------------------------
# 4) The six 24/12/2025 B1 "Session 14" rows flip from Not Recorded
#    (pink) to Recorded (green) now that attendance has been entered.
# -----------------------------------------------------------------
$newlyRecorded = @(
    @{ Row = 171; Prev = 170; Attendance = "20/23" },   # B1D1
    @{ Row = 198; Prev = 197; Attendance = "22/30" },   # B1D2
    @{ Row = 225; Prev = 224; Attendance = "19/25" },   # B1E1
    @{ Row = 252; Prev = 251; Attendance = "20/28" },   # B1E2
    @{ Row = 279; Prev = 278; Attendance = "15/26" },   # B1F1
    @{ Row = 306; Prev = 305; Attendance = "22/29" }    # B1F2
)

foreach ($item in $newlyRecorded) {
    $row = $item.Row
    $prev = $item.Prev

    # copy the "Recorded" row styling (green fill) from the row above,
    # which already carries the correct style for every column A:I
    $ws.Range("A${prev}:I${prev}").Copy()
    $ws.Range("A${row}:I${row}").PasteSpecial(-4122)

    $ws.Range("G$row").Value = "dnasr281@gmail.com"
    $ws.Range("H$row").Value = $item.Attendance
    $ws.Range("I$row").Value = "Recorded"
}
